$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.353.31"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "2.324.90"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'533.20"
$ws.Range("E5").Value = "  +2.23%  "
$ws.Range("D6").Value = "'133.12"
$ws.Range("E6").Value = "  -2.86%  "
$ws.Range("D7").Value = "'0.994"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("D9").Value = "2.351.78"
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("E12").Value = "  -2.58%  "
$ws.Range("D13").Value = "'0.345"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'23.60"
$ws.Range("E14").Value = "  -2.23%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.747.88"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "57.389.27"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E17").Value = "  -1.77%  "
$ws.Range("D18").Value = "2.334.28"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").Value = "'341.17"
$ws.Range("E19").Value = "  +3.62%  "
$ws.Range("D20").Value = "'10.47"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").Value = "'6.93"
$ws.Range("E21").Value = "  +2.49%  "
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'62.07"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("D25").Value = "'8.81"
$ws.Range("E25").Value = "  +6.96%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "'0.991"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").Value = "'1.34"
$ws.Range("E28").Value = "  +1.77%  "
$ws.Range("D29").Value = "'170.54"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("D31").Value = "0.0₃0728"
$ws.Range("E31").Value = "  -2.90%  "
$ws.Range("E32").Value = "  -2.43%  "
$ws.Range("D33").Value = "'18.56"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("E36").Value = "  -2.54%  "
$ws.Range("D37").Value = "'4.03"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").Value = "'0.910"
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("D39").Value = "'1.59"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").Value = "'39.13"
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("D41").Value = "'149.08"
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("D42").Value = "'0.378"
$ws.Range("E42").Value = "  -2.04%  "
$ws.Range("D43").Value = "'3.60"
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("D44").Value = "'282.43"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "'5.14"
$ws.Range("E45").Value = "  -3.60%  "
$ws.Range("D46").Value = "'0.0931"
$ws.Range("E46").Value = "  -1.05%  "
$ws.Range("D47").Value = "'0.0506"
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("D49").Value = "'18.64"
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("D50").Value = "'0.0218"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("E51").Value = "  -0.72%  "
